# Update "想去人数" (want-to-go count) figures that changed on re-scrape.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 84
$wsExhibit.Range("F4").Value = 487
$wsExhibit.Range("F5").Value = 4777
$wsExhibit.Range("F9").Value = 736
$wsExhibit.Range("F10").Value = 214

# Sheet "演出" (performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 23

# Sheet "全部类型" (all types combined)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 84
$wsAll.Range("F4").Value = 487
$wsAll.Range("F5").Value = 4777
$wsAll.Range("F9").Value = 736
$wsAll.Range("F10").Value = 23
$wsAll.Range("F11").Value = 214
